# "add personal info - room info"
# Updates quantity counts across several sheets, and touches column D on the
# "기타" sheet (a new room-info column was started there) leaving it blank
# for row 14 so the sheet's used range grows to include column D.

$wb = $excel.ActiveWorkbook

# 매점판매 (store sales): 대패삼겹살 quantity 26 -> 27
$ws2 = $wb.Worksheets.Item("매점판매")
$ws2.Range("C5").Value = 27

# 장의용품 (funeral supplies): 맥주 quantity 10 -> 13
$ws3 = $wb.Worksheets.Item("장의용품")
$ws3.Range("C5").Value = 13

# 상복 (mourning clothes): 맥주 quantity 13 -> 16
$ws4 = $wb.Worksheets.Item("상복")
$ws4.Range("C11").Value = 16

# 기타 (other): 치즈김밥 233 -> 231, 치킨 0 -> 1, and touch D14 (new room-info
# column) so the used range extends to column D
$ws5 = $wb.Worksheets.Item("기타")
$ws5.Range("C10").Value = 231
$ws5.Range("C13").Value = 1
$ws5.Range("D14").Value = "x"
$ws5.Range("D14").Value = ""
